$d = $word.ActiveDocument

$d.Content.Find.Execute("71×26=", $true, $false, $false, $false, $false, $true, 1, $false, "50×39=", 2) | Out-Null
$d.Content.Find.Execute("88×11=", $true, $false, $false, $false, $false, $true, 1, $false, "62×67=", 2) | Out-Null
$d.Content.Find.Execute("38×94=", $true, $false, $false, $false, $false, $true, 1, $false, "74×56=", 2) | Out-Null
$d.Content.Find.Execute("35×63=", $true, $false, $false, $false, $false, $true, 1, $false, "90×82=", 2) | Out-Null
$d.Content.Find.Execute("85×47=", $true, $false, $false, $false, $false, $true, 1, $false, "33×27=", 2) | Out-Null
$d.Content.Find.Execute("91×17=", $true, $false, $false, $false, $false, $true, 1, $false, "93×19=", 2) | Out-Null
$d.Content.Find.Execute("27×45=", $true, $false, $false, $false, $false, $true, 1, $false, "53×82=", 2) | Out-Null
$d.Content.Find.Execute("94×72=", $true, $false, $false, $false, $false, $true, 1, $false, "81×60=", 2) | Out-Null
$d.Content.Find.Execute("13×87=", $true, $false, $false, $false, $false, $true, 1, $false, "52×32=", 2) | Out-Null
$d.Content.Find.Execute("62×73=", $true, $false, $false, $false, $false, $true, 1, $false, "41×26=", 2) | Out-Null
$d.Content.Find.Execute("85×63=", $true, $false, $false, $false, $false, $true, 1, $false, "30×86=", 2) | Out-Null
$d.Content.Find.Execute("50×94=", $true, $false, $false, $false, $false, $true, 1, $false, "58×87=", 2) | Out-Null
$d.Content.Find.Execute("80×71=", $true, $false, $false, $false, $false, $true, 1, $false, "31×27=", 2) | Out-Null
$d.Content.Find.Execute("50×36=", $true, $false, $false, $false, $false, $true, 1, $false, "17×75=", 2) | Out-Null
$d.Content.Find.Execute("37×76=", $true, $false, $false, $false, $false, $true, 1, $false, "73×46=", 2) | Out-Null
$d.Content.Find.Execute("16×67=", $true, $false, $false, $false, $false, $true, 1, $false, "47×55=", 2) | Out-Null
$d.Content.Find.Execute("20×14=", $true, $false, $false, $false, $false, $true, 1, $false, "51×54=", 2) | Out-Null
$d.Content.Find.Execute("43×22=", $true, $false, $false, $false, $false, $true, 1, $false, "81×23=", 2) | Out-Null
$d.Content.Find.Execute("69×83=", $true, $false, $false, $false, $false, $true, 1, $false, "95×31=", 2) | Out-Null
$d.Content.Find.Execute("72×52=", $true, $false, $false, $false, $false, $true, 1, $false, "29×98=", 2) | Out-Null
$d.Content.Find.Execute("63×11=", $true, $false, $false, $false, $false, $true, 1, $false, "77×37=", 2) | Out-Null
$d.Content.Find.Execute("97×82=", $true, $false, $false, $false, $false, $true, 1, $false, "36×96=", 2) | Out-Null
$d.Content.Find.Execute("67×20=", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=", 2) | Out-Null
$d.Content.Find.Execute("95×77=", $true, $false, $false, $false, $false, $true, 1, $false, "70×32=", 2) | Out-Null
$d.Content.Find.Execute("14×84=", $true, $false, $false, $false, $false, $true, 1, $false, "54×92=", 2) | Out-Null
